$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '94.833.31'
$ws.Range('E2').Value = '  +3.33%  '
$ws.Range('D3').Value = '3.068.13'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.02'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '601.49'
$ws.Range('E6').Value = '  -1.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.10'
$ws.Range('E7').Value = '  +1.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.377'
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '3.064.00'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.779'
$ws.Range('E11').Value = '  +1.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.195'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').Value = '94.278.67'
$ws.Range('E13').Value = '  +2.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000235'
$ws.Range('E14').Value = '  -3.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.30'
$ws.Range('E15').Value = '  -1.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.95'
$ws.Range('E16').Value = '  -1.28%  '
$ws.Range('D17').Value = '3.628.83'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').Value = '3.079.82'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.52'
$ws.Range('E19').Value = '  -7.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.09'
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '436.62'
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.54'
$ws.Range('E22').Value = '  -3.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.67'
$ws.Range('E23').Value = '  -3.93%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000186'
$ws.Range('E24').Value = '  -3.38%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.44'
$ws.Range('E25').Value = '  -1.99%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '83.94'
$ws.Range('E26').Value = '  -1.81%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.49'
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '3.218.84'
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.131'
$ws.Range('E30').Value = '  +3.21%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.240'
$ws.Range('E31').Value = '  +2.02%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.174'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '8.82'
$ws.Range('E34').Value = '  -2.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '25.03'
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.149'
$ws.Range('E36').Value = '  -5.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.22'
$ws.Range('E37').Value = '  -6.20%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '482.65'
$ws.Range('E38').Value = '  +3.91%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.06'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('B40').Value = 'PancakeSwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.84'
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('B41').Value = 'MantraDAO'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.72'
$ws.Range('E41').Value = '  -3.83%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.426'
$ws.Range('E42').Value = '  -1.66%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.21'
$ws.Range('E43').Value = '  -4.02%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.10'
$ws.Range('E45').Value = '  -4.02%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '160.40'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.667'
$ws.Range('E47').Value = '  -1.27%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.78'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.60'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000269'
$ws.Range('E51').Value = '  +11.71%  '
